$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A186").Value = "2023-12-11 13:02:26"
$ws.Range("B186").Value = 0.0016

$ws.Range("A187").Value = "2023-12-11 13:02:58"
$ws.Range("B187").Value = 0.002

$ws.Range("A188").Value = "2023-12-11 13:03:09"
$ws.Range("B188").Value = 0.0002

$ws.Range("A189").Value = "2023-12-11 13:03:17"
$ws.Range("B189").Value = 0.0006000000000000001
